# Update loading_percent values for the 380 kV case
# Columns updated (data columns, in order): C, D, E, F, G, J, K, L, N
# Rows updated: 2..25 (24 data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","J","K","L","N")

$newValues = @(
    @(13.26178636439901,3.556611813870331,13.40283408883736,67.67035775131886,3.858313013601256,12.69443507458476,33.40599697711745,9.276437674311818,23.30874401756238),
    @(13.24484011048976,3.552213494770433,13.42275561319682,67.63492924954576,3.863566877662709,12.72424311231219,33.20914451126488,9.295975581513707,23.38382902407303),
    @(13.23738718194863,3.549571841354696,13.43711456427538,67.62935827139879,3.866955320474344,12.74432109590926,33.09673558127737,9.308781006085791,23.43204865132743),
    @(13.2350931439131,3.548510155288024,13.44350050168388,67.63114583958547,3.868377195775272,12.75294931616154,33.05308786768385,9.314203365145174,23.45223202709349),
    @(13.23475711299172,3.548334757411008,13.44459315659905,67.6316873617598,3.868615782211197,12.75440897261197,33.04597153151003,9.315116085093395,23.45561571551541),
    @(13.23735323438498,3.54955746313239,13.43719852334326,67.6293659661269,3.86697432989844,12.7444356524425,33.09613814726274,9.308853306946135,23.43231868981027),
    @(13.25533118363274,3.555083032746023,13.40926138617762,67.65477726483707,3.860090919408241,12.70434418800381,33.33638720044685,9.283006780959504,23.3341945920386),
    @(13.313955779264,3.566387707307813,13.37136706033461,67.83341424391364,3.84787396617846,12.6398259528601,33.87303760995298,9.238714732357751,23.15851964493467),
    @(13.37117135378045,3.574982678543764,13.3538391776494,68.04358957360964,3.839667626053562,12.60103721484394,34.30490319411032,9.210034059546288,23.03958494672329),
    @(13.40023844242285,3.578956554777186,13.34810686844119,68.15637478507234,3.836098911545052,12.58526503482144,34.50897054346521,9.197817053524956,22.98766315730736),
    @(13.41167871638007,3.580470624852044,13.34625849181139,68.20155094911206,3.834770975372672,12.5795621644469,34.58728875632644,9.193309537633738,22.96831444753398),
    @(13.40919564930971,3.58014413054545,13.34664223752081,68.19171181644606,3.835055929737603,12.58077837814332,34.57037601718448,9.194275036927479,22.97246763326402),
    @(13.40117098274397,3.579080932551192,13.34794834157003,68.16004204272819,3.83598919232817,12.58479044867193,34.51539324439939,9.197443839707512,22.98606506020878),
    @(13.39631193412355,3.578430899220243,13.34879034402781,68.14096451360147,3.836563892244659,12.5872830932783,34.48184887666854,9.199400277792135,22.99443460670755),
    @(13.36933259382522,3.574724267855728,13.35425890176542,68.03656421941086,3.839904142463398,12.60210568577172,34.29171564616083,9.210849123722245,23.04302199561226),
    @(13.3535575316161,3.572466788994567,13.35818776727489,67.97691584138515,3.84199525371323,12.61167885654112,34.1769877870849,9.218084803763546,23.07338708223677),
    @(13.3447704851172,3.571174365606423,13.36065850295799,67.94422465747489,3.843213489644019,12.61736136556468,34.11171885113934,9.222324729489204,23.09105772521445),
    @(13.34184463254542,3.570737805624983,13.36153128096893,67.93343374685512,3.843628628296341,12.61931562930549,34.08974503199153,9.223773735324739,23.09707601365877),
    @(13.35520720391371,3.572706478225057,13.35774769913102,67.98309817374312,3.841771050119761,12.61064152972406,34.18912666813726,9.217306468636176,23.070133410101),
    @(13.40351629841204,3.579392968119724,13.34755595942578,68.16927731060052,3.835714435263552,12.58360468341493,34.53151515238503,9.196509865710143,22.98206267970258),
    @(13.43761193158825,3.583816916825959,13.34277378461338,68.30533151025413,3.831892737836567,12.56750680200059,34.76133594460654,9.183610283110829,22.92632707587775),
    @(13.41918502366916,3.581450818262632,13.34515422571269,68.23140316927862,3.833920005865045,12.57595455863483,34.63814032173767,9.190431877642572,22.95590762052938),
    @(13.35446050764507,3.572598097810671,13.35794599367703,67.98029815289264,3.841872362626422,12.61110994839128,34.18363652803153,9.217658104608381,23.07160373067052),
    @(13.29560280236574,3.563278944101759,13.37980828459988,67.77124586476073,3.851043004866387,12.65576853329297,33.72107605740167,9.250016297268189,23.20425946182211)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}
